$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 totals
$ws.Range("B2").Value = 924
$ws.Range("C2").Value = 9199
$ws.Range("D2").Value = 41
$ws.Range("E2").Value = 9240
$ws.Range("F2").Value = 2991
$ws.Range("G2").Value = 5240

# Update row 4 values ("test lokaal" -> "pre-meeting")
$ws.Range("A4").Value = "pre-meeting"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 23
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 30
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = 17

# Delete row 5 entirely (shift cells up)
$ws.Range("A5:I5").Delete()
